$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# New data for rows 2-12 (replaces old rows 2-9 entirely)
$data = @(
    @("Tienda - Velázquez", "V2", "BAR",             45702, 8877, "Mañana", "EUROS",        488.1,   38,      7),
    @("Tienda - Velázquez", "V2", "BAR",             45702, 8877, "Mañana", "TARJETA VISA", 597.95,  100.55,  20),
    @("Tienda - Velázquez", "V1", "SERVIDOR TIENDA", 45702, 8876, "Mañana", "EUROS",        908.11,  488.06,  60),
    @("Tienda - Velázquez", "V1", "SERVIDOR TIENDA", 45702, 8876, "Mañana", "SMS",          0,       5.4,     1),
    @("Tienda - Velázquez", "V1", "SERVIDOR TIENDA", 45702, 8876, "Mañana", "TARJETA VISA", 1063.98, 1084.13, 105),
    @("Tienda - Velázquez", "V1", "SERVIDOR TIENDA", 45702, 8879, "Mañana", "EUROS",        1199.5,  830.6,   82),
    @("Tienda - Velázquez", "V1", "SERVIDOR TIENDA", 45702, 8879, "Mañana", "TARJETA VISA", 2825.69, 2791.89, 229),
    @("Tienda - Velázquez", "V2", "BAR",             45702, 8878, "Mañana", "EUROS",        646.2,   189,     22),
    @("Tienda - Velázquez", "V2", "BAR",             45702, 8878, "Mañana", "TARJETA VISA", 1226.39, 468.34,  49),
    @("Tienda - Velázquez", "V1", "SERVIDOR TIENDA", 45702, 8880, "Mañana", "EUROS",        2161.42, 685.51,  66),
    @("Tienda - Velázquez", "V1", "SERVIDOR TIENDA", 45702, 8880, "Mañana", "TARJETA VISA", 3045.82, 1528.91, 105)
)

$startRow = 2
for ($i = 0; $i -lt $data.Length; $i++) {
    $row = $startRow + $i
    $rec = $data[$i]
    $ws.Cells.Item($row, 1).Value = $rec[0]
    $ws.Cells.Item($row, 2).Value = $rec[1]
    $ws.Cells.Item($row, 3).Value = $rec[2]
    $ws.Cells.Item($row, 4).Value = $rec[3]
    $ws.Cells.Item($row, 5).Value = $rec[4]
    $ws.Cells.Item($row, 6).Value = $rec[5]
    $ws.Cells.Item($row, 7).Value = $rec[6]
    $ws.Cells.Item($row, 8).Value = $rec[7]
    $ws.Cells.Item($row, 9).Value = $rec[8]
    $ws.Cells.Item($row, 10).Value = $rec[9]

    # Keep the same per-column number formats used throughout the table
    $ws.Cells.Item($row, 4).NumberFormat = "DD/MM/YYYY"
    $ws.Cells.Item($row, 9).NumberFormat = "#,##0.00"
    $ws.Cells.Item($row, 10).NumberFormat = "#,##0"
}
